# Updates "想去人数" (want-to-go count, column F) figures across the workbook,
# inserts a newly-announced show ("上海·今泉爱夏  巡演") into the "演出" sheet,
# and appends a newly-announced event ("上海·「PLAVE with animate cafe」") to the
# "本地生活" sheet. Mirrors the upstream scraper's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) — refresh "want to go" counts
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 244
    3  = 581
    6  = 3179
    7  = 2754
    11 = 353
    12 = 290
    14 = 5723
    16 = 1024
    17 = 61
    18 = 166
    19 = 85
    21 = 1243
    23 = 10
    24 = 126
    25 = 334
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) — refresh counts, then insert the new show
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    5  = 24
    8  = 337
    9  = 61
    21 = 54
    25 = 4033
    29 = 209
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Range("F$row").Value = $showUpdates[$row]
}

# Insert a new row above the old row 33 ("夏川里美…"), pushing everything
# from row 33 onward down by one (row 33->34, 34->35, 35->36).
$wsShow.Rows.Item(33).Insert()

# Re-apply the bordered/bold/centred index-column look to the new A33 (the
# plain row insert alone leaves it unbordered) before filling in its value.
$wsShow.Range("A34").Copy() | Out-Null
$wsShow.Range("A33").PasteSpecial(-4122) | Out-Null

$wsShow.Range("A33").Value = 32
$wsShow.Range("B33").NumberFormat = "@"
$wsShow.Range("B33").Value = "2024.04.21"
$wsShow.Range("C33").Value = "上海·今泉爱夏  巡演"
$wsShow.Range("D33").Value = "瑞虹路188号3楼 Modernsky Lab"
$wsShow.Range("E33").Value = "2024.04.21 20:00-04.21 21:30"
$wsShow.Range("F33").Value = 0
$wsShow.Range("G33").Value = "不可售"
$wsShow.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=81891"
$wsShow.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202402/YJENeaUi1708313389899.jpeg"

# The rows that used to be 33/34/35 are now 34/35/36; the refreshed data
# bumps their index counter (column A) up by one from the plain row-shift
# value, and the last one's want-to-go count also increments.
$wsShow.Range("A34").Value = 33
$wsShow.Range("A35").Value = 34
$wsShow.Range("A36").Value = 35
$wsShow.Range("F36").Value = 10

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life) — refresh counts, then append the new event
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$localUpdates = @{
    5  = 2567
    6  = 1116
    9  = 1444
    10 = 402
}
foreach ($row in $localUpdates.Keys) {
    $wsLocal.Range("F$row").Value = $localUpdates[$row]
}

# Append a brand-new row 13 ("PLAVE with animate cafe"), copying the
# look-and-feel (borders/bold/alignment) of the existing index column.
$wsLocal.Range("A12").Copy() | Out-Null
$wsLocal.Range("A13").PasteSpecial(-4122) | Out-Null
$wsLocal.Range("A13").Value = 12

$wsLocal.Range("B13").NumberFormat = "@"
$wsLocal.Range("B13").Value = "2024.03.01"
$wsLocal.Range("C13").Value = "上海·「PLAVE with animate cafe」"
$wsLocal.Range("D13").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$wsLocal.Range("E13").Value = "2024.03.01 00:00-03.25 23:59"
$wsLocal.Range("F13").Value = 237
$wsLocal.Range("G13").Value = 30
$wsLocal.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=81873"
$wsLocal.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202402/7QENUAuN1708247451105.png"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types) — refresh counts (same figures as above sheets)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    5  = 2567
    6  = 1116
    7  = 1444
    8  = 402
    11 = 244
    12 = 581
    13 = 24
    14 = 3179
    15 = 2754
    19 = 353
    20 = 337
    21 = 61
    22 = 290
    26 = 1024
    28 = 61
    29 = 166
    30 = 85
    36 = 54
    39 = 1243
    41 = 209
    43 = 10
    46 = 126
    47 = 334
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
